# Append new job listing row (2026-02-13 13:08 JST scrape) into the
# "ランサーズ" sheet, shifting the existing rows down, and refresh every
# row's "取得日時" (scraped-at) timestamp to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-02-13 13:08:42"

# --- 1. Insert a fresh row right above the old row 5 (the "AI agent
#        partner" posting), pushing rows 5-12 down to 6-13. ------------
$ws.Rows.Item(5).Insert()

# --- 2. Populate the newly inserted row 5 with the new posting. -------
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5460294"
$ws.Range("G5").Value = 375
$ws.Range("H5").Value = "🔥AI,Ai ◆開発"

# --- 3. Refresh the "取得日時" timestamp for every data row (2-13), -----
#        including the ones that merely shifted down. ------------------
$lastRow = 13
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 4. Rebuild the URL hyperlinks in column F (row insert does not ---
#        shift the existing hyperlink anchors in this engine, so drop
#        them all and re-add in row order, which also re-applies the
#        "Hyperlink" cell style). ---------------------------------------
$ws.Range("F2").Hyperlinks.Delete()

$urls = @{
    2  = "https://www.lancers.jp/work/detail/5491124"
    3  = "https://www.lancers.jp/work/detail/5490911"
    4  = "https://www.lancers.jp/work/detail/5450864"
    5  = "https://www.lancers.jp/work/detail/5460294"
    6  = "https://www.lancers.jp/work/detail/5490828"
    7  = "https://www.lancers.jp/work/detail/5491190"
    8  = "https://www.lancers.jp/work/detail/5473940"
    9  = "https://www.lancers.jp/work/detail/5490679"
    10 = "https://www.lancers.jp/work/detail/5491203"
    11 = "https://www.lancers.jp/work/detail/5477871"
    12 = "https://www.lancers.jp/work/detail/5491086"
    13 = "https://www.lancers.jp/work/detail/5490905"
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = $urls[$r]
    $ws.Hyperlinks.Add($cell, $urls[$r])
}
